# "commit from main to lib" — the workbook's lib copy gains a literal
# value in A1 (the original "main" copy had an empty sheet).
#
# Note: the diff also shows the mc:AlternateContent/x15ac:absPath url
# changing (Mac Excel's "file was last seen at this folder" breadcrumb,
# rewritten to the new /sub/ location on save). That attribute isn't
# backed by any property on the Application/Workbook/Worksheet/Range
# object model (no Workbook.AbsPath, and Path/FullName/SaveAs don't
# round-trip into it here) — it's host/save-path metadata outside what
# COM automation can drive, so it's intentionally left alone rather than
# faked via an unrelated call.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 123123
